$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I3").Value = "pink"
$ws.Range("H1").Value = "App Weather"
$ws.Range("I1").Value = "App color"
$ws.Range("J1").Value = "App Hue Code"

$ws.Range("J2").Value = 25
$ws.Range("J3").Value = 83
$ws.Range("J4").Value = 75
$ws.Range("J5").Value = 10
$ws.Range("J6").Value = 70

$ws.Range("A1").Font.Bold = $true
$ws.Range("E1").Font.Bold = $true
$ws.Range("F1").Font.Bold = $true
$ws.Range("H1").Font.Bold = $true
$ws.Range("I1").Font.Bold = $true
$ws.Range("J1").Font.Bold = $true

$ws.Columns.Item(8).ColumnWidth = 12.666666666666666
$ws.Columns.Item(9).ColumnWidth = 9.333333333333334
$ws.Columns.Item(10).ColumnWidth = 13.833333333333334

$ws.Range("J3").Select()
